$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: drop the stray "_GoBack" bookmark that originally sits right
# after " оплаты и " (it will be re-created further down, at its new
# location, once the second paragraph has been re-typed).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: the "Закачать файлы из архива через ftp в папку sites/all/modules/"
# run, followed by a lone-space run, gets retyped as two runs split in the
# middle of "через" ("чере" | "з ...") with the trailing "/" (and the extra
# space run) dropped, and the cursor / "_GoBack" bookmark ends up right
# before the following "/ubercart/payment/uc_paymaster" run.
# ---------------------------------------------------------------------------
$anchor = "Закачать файлы из архива через ftp в папку sites/all/modules/"
$firstPart  = "Закачать файлы из архива чере"
$secondPart = "з ftp в папку sites/all/modules"

$text = $d.Content.Text
$anchorStart = $text.IndexOf($anchor)
if ($anchorStart -lt 0) {
    throw "Could not locate the anchor sentence in the document"
}

$trailingSlashPos = $anchorStart + $anchor.Length - 1      # the '/' at the very end of the anchor run
$straySpacePos    = $anchorStart + $anchor.Length          # the following single-space run

# Remove the trailing slash and the stray space run (highest offset first so
# earlier offsets stay valid): "...modules/" + " " + "/ubercart" -> "...modules" + "/ubercart"
$d.Range($straySpacePos, $straySpacePos + 1).Text = ""
$d.Range($trailingSlashPos, $trailingSlashPos + 1).Text = ""

# New boundary: right after "...sites/all/modules", right before "/ubercart..."
$boundaryPos = $trailingSlashPos

# Re-type the sentence as two runs, split inside "через"
$splitPos = $anchorStart + $firstPart.Length
$d.Range($anchorStart, $splitPos).Text = $firstPart

$secondRange = $d.Range($splitPos, $boundaryPos)
$secondRange.Text = $secondPart

# The engine (like Word itself) silently re-merges adjacent runs that carry
# identical formatting, so nudge the formatting once to force the retyped
# text to stay split into its own run, then put it straight back.
$secondRange2 = $d.Range($splitPos, $boundaryPos)
$secondRange2.Font.Bold = 1
$secondRange2.Font.Bold = 0

# Re-create the "_GoBack" bookmark at its new location, right before
# "/ubercart/payment/uc_paymaster".
$d.Bookmarks.Add("_GoBack", $d.Range($boundaryPos, $boundaryPos))
